# Apply "MAJ mapping suite review de NRISS" edits to the workbook.

$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata": swap Name/Title values and bump the Date ---
$wsMeta = $wb.Worksheets.Item("Metadata")

# Name row no longer carries the old "Mapping Métier/CDA/FHIR..." text
$wsMeta.Cells.Item(4, 2).Value2 = $null

# Title row now carries what used to be the Name's value
$wsMeta.Cells.Item(5, 2).Value2 = 'Mapping Métier/CDA/FHIR : "Utilisation de dispositif médical"'

# Date bump
$wsMeta.Cells.Item(8, 2).Value2 = "2026-01-07T15:20:53+00:00"

# --- Sheet "Mapping Table 0": dotted entryRelationship targets -> colon form ---
$wsT0 = $wb.Worksheets.Item("Mapping Table 0")

$wsT0.Cells.Item(8, 4).Value2 = "FRCDADispositifMedical.entryRelationship:frEnRapportAvecALD"
$wsT0.Cells.Item(9, 4).Value2 = "FRCDADispositifMedical.entryRelationship:frEnRapportAvecAccidentTravail"
$wsT0.Cells.Item(10, 4).Value2 = "FRCDADispositifMedical.entryRelationship:frEnRapportAvecPrevention"

# --- Sheet "Mapping Table 1": dotted extension/reasonReference/entryRelationship targets -> colon form ---
$wsT1 = $wb.Worksheets.Item("Mapping Table 1")

$wsT1.Cells.Item(6, 4).Value2 = "FRDeviceUseStatementDocument.source.extension:performer"

$wsT1.Cells.Item(8, 1).Value2 = "FRCDADispositifMedical.entryRelationship:frEnRapportAvecALD"
$wsT1.Cells.Item(8, 4).Value2 = "FRDeviceUseStatementDocument.reasonReference:EnRapportAvecALD"

$wsT1.Cells.Item(9, 1).Value2 = "FRCDADispositifMedical.entryRelationship:frEnRapportAvecAccidentTravail"
$wsT1.Cells.Item(9, 4).Value2 = "FRDeviceUseStatementDocument.reasonReference:EnRapportAvecAccidentTravail"

$wsT1.Cells.Item(10, 1).Value2 = "FRCDADispositifMedical.entryRelationship:frEnRapportAvecPrevention"
$wsT1.Cells.Item(10, 4).Value2 = "FRDeviceUseStatementDocument.reasonReference:EnRapportAvecLaPrevention"

$wsT1.Cells.Item(11, 1).Value2 = "FRCDADispositifMedical.entryRelationship:frNonRemboursable"
$wsT1.Cells.Item(11, 4).Value2 = "FRDeviceUseStatementDocument.extension:notCovered"
